# Slide 7 ("Conclusiones y acuerdos"), content placeholder shape:
#   - change paragraph 1 run language es-ES -> es-MX and drop its endParaRPr
#   - extend paragraph 2's sentence
#   - append three new bullet paragraphs
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# 1) Normalize language to es-MX on the whole range while it is still the
#    original two paragraphs (updates both existing runs' rPr/endParaRPr).
$tr.LanguageID = "es-MX"

$t1 = "Realizar un estudio sobre el estado del arte y soluciones a problemas similares."

# 2) Drop the stray <a:endParaRPr> left on paragraph 1: delete paragraph 1
#    (text + its paragraph mark) and re-insert the same text + mark at the
#    very start, which regenerates the paragraph without an endParaRPr.
$para1incl = $tr.Characters(1, $t1.Length + 1)
[void]$para1incl.Delete()
[void]$tr.InsertBefore($t1 + "`r")

# 3) Paragraph 2: "Realizar un prototipo de solución para mostrar al
#    cliente." -> "...al cliente en la siguiente visita."
$fullLen = $tr.Length
$para2Start = $t1.Length + 2
$para2TotalLen = $fullLen - $para2Start + 1
$para2 = $tr.Characters($para2Start, $para2TotalLen)
$para2.Text = "Realizar un prototipo de solución para mostrar al cliente en la siguiente visita."

# 4) Add three new paragraphs after the (now updated) second paragraph.
$t3 = "Siguiente fecha oficial establecida en la tercera semana de septiembre."
$t4 = "Por parte del contacto se acuerda que es posible buscarlo para resolver dudas si es necesario antes de la siguiente fecha pactada."
$t5 = "Contacto acordó enviar al equipo de trabajo fotografías para documentar la zona."
[void]$tr.InsertAfter("`r" + $t3 + "`r" + $t4 + "`r" + $t5)
